$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 1107.2307
$ws.Range("J41").Value = 1499.75
$ws.Range("L41").Value = 1499.75
$ws.Range("N41").Value = -2379.75
$ws.Range("H100").Value = 1681.9333
$ws.Range("I100").Value = 1702.2307
$ws.Range("K100").Value = 1702.2307
$ws.Range("M100").Value = -1161.2307
$ws.Range("H104").Value = 467.16666
$ws.Range("I104").Value = 420.6
$ws.Range("K104").Value = 1261.8
$ws.Range("M104").Value = 485.1999999999998
$ws.Range("H112").Value = 59231.61
$ws.Range("I112").Value = 2250
$ws.Range("K112").Value = 6750
$ws.Range("M112").Value = -5642
$ws.Range("H113").Value = 2427.3333
$ws.Range("I113").Value = 2392.077
$ws.Range("K113").Value = 2392.077
$ws.Range("M113").Value = 861.9229999999998
$ws.Range("H116").Value = 5495.3335
$ws.Range("I116").Value = 5510.3
$ws.Range("K116").Value = 5510.3
$ws.Range("M116").Value = -2068.3
$ws.Range("H138").Value = 5036.49
$ws.Range("I138").Value = 13951.25
$ws.Range("K138").Value = 41853.75
$ws.Range("M138").Value = -36713.75
$ws.Range("H141").Value = 1897
$ws.Range("I141").Value = 1897
$ws.Range("K141").Value = 5691
$ws.Range("M141").Value = -511

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 26924.666
$ws.Range("I25").Value = 1800
$ws.Range("K25").Value = 1800
$ws.Range("M25").Value = -1398
$ws.Range("H32").Value = 2642.1667
$ws.Range("I32").Value = 1086.5385
$ws.Range("K32").Value = 1086.5385
$ws.Range("M32").Value = -799.5385000000001
$ws.Range("H45").Value = 47731.953
$ws.Range("I45").Value = 68340.2
$ws.Range("K45").Value = 68340.2
$ws.Range("M45").Value = -67963.2
$ws.Range("H61").Value = 4064692.8
$ws.Range("I61").Value = 94780.5
$ws.Range("K61").Value = 94780.5
$ws.Range("M61").Value = -94568.5
$ws.Range("H74").Value = 638140.9399999999
$ws.Range("I74").Value = 1083.2941
$ws.Range("K74").Value = 1083.2941
$ws.Range("M74").Value = -209.2941000000001
$ws.Range("H77").Value = 638140.9399999999
$ws.Range("I77").Value = 1083.2941
$ws.Range("K77").Value = 5416.4705
$ws.Range("M77").Value = -1048.4705
$ws.Range("H94").Value = 45000
$ws.Range("I94").Value = 20000
$ws.Range("K94").Value = 20000
$ws.Range("M94").Value = -19099
$ws.Range("H110").Value = 1346.65
$ws.Range("I110").Value = 1385.5
$ws.Range("K110").Value = 1385.5
$ws.Range("M110").Value = 659.5
$ws.Range("H122").Value = 1632.5333
$ws.Range("I122").Value = 1540.6666
$ws.Range("K122").Value = 4621.9998
$ws.Range("M122").Value = -2171.9998
$ws.Range("H132").Value = 2956.3044
$ws.Range("I132").Value = 2979
$ws.Range("J132").Value = 2921
$ws.Range("K132").Value = 8937
$ws.Range("L132").Value = 8763
$ws.Range("M132").Value = -6407
$ws.Range("N132").Value = -13823
$ws.Range("H136").Value = 4064692.8
$ws.Range("I136").Value = 94780.5
$ws.Range("K136").Value = 284341.5
$ws.Range("M136").Value = -281791.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 944.5333000000001
$ws.Range("J20").Value = 914.5714
$ws.Range("L20").Value = 914.5714
$ws.Range("N20").Value = -1408.5714

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H38").Value = 11950
$ws.Range("J38").Value = 11950
$ws.Range("L38").Value = 11950
$ws.Range("N38").Value = -12704
$ws.Range("H46").Value = 11950
$ws.Range("J46").Value = 11950
$ws.Range("L46").Value = 11950
$ws.Range("N46").Value = -12372
$ws.Range("H86").Value = 13948.272
$ws.Range("I86").Value = 4705.2856
$ws.Range("J86").Value = 30123.5
$ws.Range("K86").Value = 4705.2856
$ws.Range("L86").Value = 30123.5
$ws.Range("M86").Value = -3582.2856
$ws.Range("N86").Value = -32369.5
$ws.Range("H89").Value = 13948.272
$ws.Range("I89").Value = 4705.2856
$ws.Range("J89").Value = 30123.5
$ws.Range("K89").Value = 23526.428
$ws.Range("L89").Value = 150617.5
$ws.Range("M89").Value = -17910.428
$ws.Range("N89").Value = -161849.5
$ws.Range("H92").Value = 36700.332
$ws.Range("J92").Value = 36700.332
$ws.Range("L92").Value = 36700.332
$ws.Range("N92").Value = -41692.332
$ws.Range("H96").Value = 15718
$ws.Range("J96").Value = 15718
$ws.Range("L96").Value = 15718
$ws.Range("N96").Value = -21210
$ws.Range("H99").Value = 57780484
$ws.Range("I99").Value = 5002609
$ws.Range("J99").Value = 100002780
$ws.Range("K99").Value = 5002609
$ws.Range("L99").Value = 100002780
$ws.Range("M99").Value = -5001111
$ws.Range("N99").Value = -100005776
$ws.Range("H126").Value = 57780484
$ws.Range("I126").Value = 5002609
$ws.Range("J126").Value = 100002780
$ws.Range("K126").Value = 15007827
$ws.Range("L126").Value = 300008340
$ws.Range("M126").Value = -15005357
$ws.Range("N126").Value = -300013280
$ws.Range("H132").Value = 25722036
$ws.Range("I132").Value = 129669.125
$ws.Range("K132").Value = 389007.375
$ws.Range("M132").Value = -386477.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 32462000
$ws.Range("J9").Value = 32462000
$ws.Range("L9").Value = 97386000
$ws.Range("N9").Value = -97386448
$ws.Range("H56").Value = 7249.5
$ws.Range("I56").Value = 7249.5
$ws.Range("K56").Value = 7249.5
$ws.Range("M56").Value = -6719.5
$ws.Range("H87").Value = 16794.75
$ws.Range("I87").Value = 3166
$ws.Range("J87").Value = 24972
$ws.Range("K87").Value = 9498
$ws.Range("L87").Value = 74916
$ws.Range("M87").Value = -8250
$ws.Range("N87").Value = -77412
$ws.Range("H90").Value = 16794.75
$ws.Range("I90").Value = 3166
$ws.Range("J90").Value = 24972
$ws.Range("K90").Value = 28494
$ws.Range("L90").Value = 224748
$ws.Range("M90").Value = -22254
$ws.Range("N90").Value = -237228
$ws.Range("H94").Value = 14679.429
$ws.Range("I94").Value = 5385
$ws.Range("K94").Value = 16155
$ws.Range("M94").Value = -15479
$ws.Range("H99").Value = 16666.818
$ws.Range("J99").Value = 19926.111
$ws.Range("L99").Value = 59778.333
$ws.Range("N99").Value = -64270.333
$ws.Range("H105").Value = 16644.857
$ws.Range("J105").Value = 17752.334
$ws.Range("L105").Value = 53257.00199999999
$ws.Range("N105").Value = -58499.00199999999
$ws.Range("H139").Value = 10419698
$ws.Range("I139").Value = 15626590
$ws.Range("K139").Value = 46879770
$ws.Range("M139").Value = -46874630

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5967.8335
$ws.Range("I70").Value = 5999.6665
$ws.Range("K70").Value = 5999.6665
$ws.Range("M70").Value = -5729.6665
$ws.Range("H73").Value = 5967.8335
$ws.Range("I73").Value = 5999.6665
$ws.Range("K73").Value = 5999.6665
$ws.Range("M73").Value = -5063.6665
$ws.Range("H113").Value = 1228.1428
$ws.Range("I113").Value = 1339.6
$ws.Range("K113").Value = 1339.6
$ws.Range("M113").Value = 830.4000000000001
$ws.Range("H126").Value = 4247
$ws.Range("I126").Value = 3846.0908
$ws.Range("K126").Value = 11538.2724
$ws.Range("M126").Value = -9068.2724
$ws.Range("H140").Value = 92930.336
$ws.Range("J140").Value = 97574.60000000001
$ws.Range("L140").Value = 97574.60000000001
$ws.Range("N140").Value = -107934.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8699.75
$ws.Range("I7").Value = 4553.6665
$ws.Range("K7").Value = 4553.6665
$ws.Range("M7").Value = -4441.6665
$ws.Range("H61").Value = 3741.9443
$ws.Range("I61").Value = 3783.8462
$ws.Range("K61").Value = 3783.8462
$ws.Range("M61").Value = -3581.8462
$ws.Range("H94").Value = 40000
$ws.Range("J94").Value = 40000
$ws.Range("L94").Value = 40000
$ws.Range("N94").Value = -41352
$ws.Range("H113").Value = 3741.9443
$ws.Range("I113").Value = 3783.8462
$ws.Range("K113").Value = 3783.8462
$ws.Range("M113").Value = -1613.8462
$ws.Range("H122").Value = 3962.6667
$ws.Range("I122").Value = 3962.6667
$ws.Range("K122").Value = 11888.0001
$ws.Range("M122").Value = -9438.000100000001
$ws.Range("H126").Value = 8699.75
$ws.Range("I126").Value = 4553.6665
$ws.Range("K126").Value = 13660.9995
$ws.Range("M126").Value = -11190.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H94").Value = 25000
$ws.Range("J94").Value = 25000
$ws.Range("L94").Value = 25000
$ws.Range("N94").Value = -26802
$ws.Range("H136").Value = 31528.242
$ws.Range("I136").Value = 42447.918
$ws.Range("J136").Value = 2409.111
$ws.Range("K136").Value = 127343.754
$ws.Range("L136").Value = 7227.333
$ws.Range("M136").Value = -124793.754
$ws.Range("N136").Value = -12327.333
